$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) values per row to reflect latest crypto data
$ws.Range("D2").Value = "30.007.00"
$ws.Range("E2").Value = "  -0.74%  "

$ws.Range("D3").Value = "1.917.38"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.17"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.83%  "

$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5037"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4022"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08244"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.111"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.08"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.64"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("D13").Value = "1.915.92"
$ws.Range("E13").Value = "  -0.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.401"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.299"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.16"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06490"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.12%  "

$ws.Range("E20").Value = "  -2.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.942"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.40%  "

$ws.Range("D23").Value = "30.053.46"
$ws.Range("E23").Value = "  -0.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.28"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.61%  "

$ws.Range("E25").Value = "  -1.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "22.28"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.31%  "

$ws.Range("D27").Value = "2.137.11"
$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.79"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.312"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.96"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.136"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1041"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.008"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.814"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02448"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.365"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06425"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.907"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.91%  "

$ws.Range("E39").Value = "  -2.43%  "

$ws.Range("E40").Value = "  -2.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6420"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.37"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.216"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.45%  "

$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.35"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6007"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.163"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.630"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.02"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.38%  "

$ws.Range("E50").Value = "  -2.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.75"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.06%  "
